$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "pozorovať súhvezdie Súhvezdie Bootes",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "pozorovať Súhvezdie Bootes",
    2
)
